$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued cells (coin name / link) - plain assignment
$textUpdates = @{
    "B6" = 'FTXToken'
    "C6" = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    "B7" = 'MXToken'
    "C7" = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    "B8" = 'BTSEToken'
    "C8" = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    "B9" = 'LiechtensteinCryptoassetsExchange'
    "C9" = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    "B10" = 'WazirX'
    "C10" = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    "B11" = 'MandalaExchangeToken'
    "C11" = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    "B12" = 'BitrueCoin'
    "C12" = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    "B13" = 'BitMartToken'
    "C13" = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    "B14" = 'BitForexToken'
    "C14" = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    "B15" = 'TigerCash'
    "C15" = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    "B16" = 'LEO'
    "C16" = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    "B17" = 'GateToken'
    "C17" = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
}

# Numeric-looking cells (price / volume%) - force text format to preserve exact string
$numericUpdates = @{
    "D2" = '309.97'
    "E2" = '1.20%'
    "D3" = '41.11'
    "E3" = '1.94%'
    "D4" = '5.121'
    "E4" = '0.15%'
    "D5" = '0.07682'
    "E5" = '1.21%'
    "D6" = '1.625'
    "E6" = '0.85%'
    "D7" = '0.9213'
    "E7" = '1.53%'
    "D8" = '2.468'
    "E8" = '1.86%'
    "D9" = '0.1227'
    "E9" = '21.31%'
    "D10" = '0.1822'
    "E10" = '3.91%'
    "D11" = '0.09161'
    "E11" = '0.05%'
    "D12" = '0.04329'
    "E12" = '3.66%'
    "D13" = '0.1051'
    "E13" = '-0.54%'
    "D14" = '0.001228'
    "E14" = '-1.43%'
    "D15" = '0.005827'
    "E15" = '-0.95%'
    "D16" = '3.353'
    "E16" = '-0.01%'
    "D17" = '4.283'
    "E17" = '0.28%'
    "D19" = '6.898'
    "E19" = '3.64%'
    "D20" = '0.1386'
    "E20" = '2.11%'
    "E21" = '-1.92%'
    "D22" = '0.04036'
    "E22" = '-3.43%'
    "D23" = '0.001263'
    "E23" = '2.93%'
    "D24" = '0.004089'
    "E24" = '0.70%'
    "D25" = '0.0001269'
    "E25" = '-2.49%'
    "E26" = '24.61%'
    "D38" = '0.02464'
    "E38" = '3.24%'
    "D39" = '0.05259'
    "E39" = '2.02%'
    "D40" = '0.007834'
    "E40" = '0.69%'
    "D41" = '0.1314'
    "E41" = '1.34%'
    "E42" = '-4.01%'
    "D43" = '0.001842'
    "E43" = '-4.08%'
    "D44" = '0.008197'
    "E44" = '-2.88%'
    "D45" = '0.3093'
    "E45" = '-6.39%'
    "D46" = '0.00006793'
    "E46" = '6.91%'
    "D47" = '0.00000000749'
    "E47" = '-0.20%'
    "D48" = '0.2251'
    "E48" = '3,068.66%'
    "D49" = '0.004087'
    "D50" = '0.00002098'
    "E50" = '-0.20%'
    "D51" = '0.0001998'
    "E51" = '-0.20%'
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

foreach ($ref in $numericUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericUpdates[$ref]
}

Write-Host "Applied $($textUpdates.Count + $numericUpdates.Count) cell updates"
